$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 8.011001628657882
$ws.Cells.Item(2, 3).Value = 6.755142719956887
$ws.Cells.Item(2, 5).Value = 13.14754306142227
$ws.Cells.Item(2, 6).Value = 16.86991607391245
$ws.Cells.Item(2, 7).Value = 23.8715216965779
$ws.Cells.Item(2, 8).Value = 13.23310748839243
$ws.Cells.Item(2, 11).Value = 8.629603700096226
$ws.Cells.Item(2, 13).Value = 12.95800125640976
$ws.Cells.Item(2, 15).Value = 19.46139881450973
$ws.Cells.Item(3, 2).Value = 7.672958716067894
$ws.Cells.Item(3, 3).Value = 6.692283743844301
$ws.Cells.Item(3, 5).Value = 12.95316928396511
$ws.Cells.Item(3, 6).Value = 15.89584955866815
$ws.Cells.Item(3, 7).Value = 24.04615748948496
$ws.Cells.Item(3, 8).Value = 13.29217416085658
$ws.Cells.Item(3, 11).Value = 8.387501410138352
$ws.Cells.Item(3, 13).Value = 12.7654934613102
$ws.Cells.Item(3, 15).Value = 19.57202063186291
$ws.Cells.Item(4, 2).Value = 7.458085905726442
$ws.Cells.Item(4, 3).Value = 6.653560038342963
$ws.Cells.Item(4, 5).Value = 12.83725739950132
$ws.Cells.Item(4, 6).Value = 15.26997757108489
$ws.Cells.Item(4, 7).Value = 24.16233712943177
$ws.Cells.Item(4, 8).Value = 13.33059209240318
$ws.Cells.Item(4, 11).Value = 8.233978172207777
$ws.Cells.Item(4, 13).Value = 12.6485670906986
$ws.Cells.Item(4, 15).Value = 19.64437482887959
$ws.Cells.Item(5, 2).Value = 7.368803943409917
$ws.Cells.Item(5, 3).Value = 6.637758667061432
$ws.Cells.Item(5, 5).Value = 12.79095244566524
$ws.Cells.Item(5, 6).Value = 15.008197319934
$ws.Cells.Item(5, 7).Value = 24.21192068588886
$ws.Cells.Item(5, 8).Value = 13.34678915737474
$ws.Cells.Item(5, 11).Value = 8.170251199112888
$ws.Cells.Item(5, 13).Value = 12.60129813252724
$ws.Cells.Item(5, 15).Value = 19.67497352588882
$ws.Cells.Item(6, 2).Value = 7.353878681304884
$ws.Cells.Item(6, 3).Value = 6.63513390672716
$ws.Cells.Item(6, 5).Value = 12.78332161821348
$ws.Cells.Item(6, 6).Value = 14.96433081551589
$ws.Cells.Item(6, 7).Value = 24.22028891053557
$ws.Cells.Item(6, 8).Value = 13.34951138826682
$ws.Cells.Item(6, 11).Value = 8.159600868551745
$ws.Cells.Item(6, 13).Value = 12.59347379710158
$ws.Cells.Item(6, 15).Value = 19.68012166688661
$ws.Cells.Item(7, 2).Value = 7.456888611132573
$ws.Cells.Item(7, 3).Value = 6.653347006560695
$ws.Cells.Item(7, 5).Value = 12.83662906259368
$ws.Cells.Item(7, 6).Value = 15.26647399323133
$ws.Cells.Item(7, 7).Value = 24.16299677882299
$ws.Cells.Item(7, 8).Value = 13.33080833855965
$ws.Cells.Item(7, 11).Value = 8.233123362218031
$ws.Cells.Item(7, 13).Value = 12.647927990695
$ws.Cells.Item(7, 15).Value = 19.64478298483609
$ws.Cells.Item(8, 2).Value = 7.896029162537027
$ws.Cells.Item(8, 3).Value = 6.73350214751334
$ws.Cells.Item(8, 5).Value = 13.07985066402796
$ws.Cells.Item(8, 6).Value = 16.5399640634477
$ws.Cells.Item(8, 7).Value = 23.92987089768259
$ws.Cells.Item(8, 8).Value = 13.25302766006859
$ws.Cells.Item(8, 11).Value = 8.547172355139763
$ws.Cells.Item(8, 13).Value = 12.89139563573054
$ws.Cells.Item(8, 15).Value = 19.49862094813845
$ws.Cells.Item(9, 2).Value = 8.694664400355242
$ws.Cells.Item(9, 3).Value = 6.889186830495156
$ws.Cells.Item(9, 5).Value = 13.58105010518545
$ws.Cells.Item(9, 6).Value = 19.00274580682531
$ws.Cells.Item(9, 7).Value = 23.54427762380483
$ws.Cells.Item(9, 8).Value = 13.11753688120166
$ws.Cells.Item(9, 11).Value = 9.121983643947454
$ws.Cells.Item(9, 13).Value = 13.37634508814354
$ws.Cells.Item(9, 15).Value = 19.24719406423778
$ws.Cells.Item(10, 2).Value = 9.238242901643519
$ws.Cells.Item(10, 3).Value = 7.002012406740696
$ws.Cells.Item(10, 5).Value = 13.95996491341444
$ws.Cells.Item(10, 6).Value = 20.67494806633232
$ws.Cells.Item(10, 7).Value = 23.30534443816785
$ws.Cells.Item(10, 8).Value = 13.02833820372544
$ws.Cells.Item(10, 11).Value = 9.51642522053937
$ws.Cells.Item(10, 13).Value = 13.73376203371653
$ws.Cells.Item(10, 15).Value = 19.08397080824636
$ws.Cells.Item(11, 2).Value = 9.475323631505415
$ws.Cells.Item(11, 3).Value = 7.052867561283764
$ws.Cells.Item(11, 5).Value = 14.13379342251054
$ws.Cells.Item(11, 6).Value = 21.3917225636224
$ws.Cells.Item(11, 7).Value = 23.20644034478092
$ws.Cells.Item(11, 8).Value = 12.98999797312735
$ws.Cells.Item(11, 11).Value = 9.689302830863248
$ws.Cells.Item(11, 13).Value = 13.8958911287724
$ws.Cells.Item(11, 15).Value = 19.01439568468466
$ws.Cells.Item(12, 2).Value = 9.563577943423891
$ws.Cells.Item(12, 3).Value = 7.07204643670733
$ws.Cells.Item(12, 5).Value = 14.1997562087586
$ws.Cells.Item(12, 6).Value = 21.65686569030329
$ws.Cells.Item(12, 7).Value = 23.17040890214407
$ws.Cells.Item(12, 8).Value = 12.97580063875364
$ws.Cells.Item(12, 11).Value = 9.753787272804272
$ws.Cells.Item(12, 13).Value = 13.95716131366872
$ws.Cells.Item(12, 15).Value = 18.98872297171405
$ws.Cells.Item(13, 2).Value = 9.544639347756661
$ws.Cells.Item(13, 3).Value = 7.067919608268963
$ws.Cells.Item(13, 5).Value = 14.18554494569816
$ws.Cells.Item(13, 6).Value = 21.60004134736742
$ws.Cells.Item(13, 7).Value = 23.17810547876824
$ws.Cells.Item(13, 8).Value = 12.97884400514739
$ws.Cells.Item(13, 11).Value = 9.739943565065971
$ws.Cells.Item(13, 13).Value = 13.94397212939678
$ws.Cells.Item(13, 15).Value = 18.99422205453228
$ws.Cells.Item(14, 2).Value = 9.482615185029076
$ws.Cells.Item(14, 3).Value = 7.05444706358527
$ws.Cells.Item(14, 5).Value = 14.1392178670244
$ws.Cells.Item(14, 6).Value = 21.4136618050453
$ws.Cells.Item(14, 7).Value = 23.20344744950519
$ws.Cells.Item(14, 8).Value = 12.9888235129188
$ws.Cells.Item(14, 11).Value = 9.694627859060931
$ws.Cells.Item(14, 13).Value = 13.90093468689025
$ws.Cells.Item(14, 15).Value = 19.01227006088646
$ws.Cells.Item(15, 2).Value = 9.444423715034391
$ws.Cells.Item(15, 3).Value = 7.046184137851368
$ws.Cells.Item(15, 5).Value = 14.11085701942395
$ws.Cells.Item(15, 6).Value = 21.29868154950795
$ws.Cells.Item(15, 7).Value = 23.21915565986136
$ws.Cells.Item(15, 8).Value = 12.99497807947602
$ws.Cells.Item(15, 11).Value = 9.666741928608255
$ws.Cells.Item(15, 13).Value = 13.87455509264836
$ws.Cells.Item(15, 15).Value = 19.02341278980397
$ws.Cells.Item(16, 2).Value = 9.222538887424205
$ws.Cells.Item(16, 3).Value = 6.998678477764781
$ws.Cells.Item(16, 5).Value = 13.94862796334798
$ws.Cells.Item(16, 6).Value = 20.62722412089977
$ws.Cells.Item(16, 7).Value = 23.3120061965673
$ws.Cells.Item(16, 8).Value = 13.03088880503152
$ws.Cells.Item(16, 11).Value = 9.504992015480267
$ws.Cells.Item(16, 13).Value = 13.72315220688851
$ws.Cells.Item(16, 15).Value = 19.08861192299438
$ws.Cells.Item(17, 2).Value = 9.083765005074232
$ws.Cells.Item(17, 3).Value = 6.969406858776045
$ws.Cells.Item(17, 5).Value = 13.84942745564351
$ws.Cells.Item(17, 6).Value = 20.20408069597325
$ws.Cells.Item(17, 7).Value = 23.37148375874567
$ws.Cells.Item(17, 8).Value = 13.05349148235706
$ws.Cells.Item(17, 11).Value = 9.404056365673355
$ws.Cells.Item(17, 13).Value = 13.63011118226577
$ws.Cells.Item(17, 15).Value = 19.12980813067812
$ws.Cells.Item(18, 2).Value = 9.00298904813101
$ws.Cells.Item(18, 3).Value = 6.95252704587834
$ws.Cells.Item(18, 5).Value = 13.79251142607727
$ws.Cells.Item(18, 6).Value = 19.95656407809801
$ws.Cells.Item(18, 7).Value = 23.40661411596621
$ws.Cells.Item(18, 8).Value = 13.06670249245424
$ws.Cells.Item(18, 11).Value = 9.345385779199674
$ws.Cells.Item(18, 13).Value = 13.57655635945497
$ws.Cells.Item(18, 15).Value = 19.15394312908168
$ws.Cells.Item(19, 2).Value = 8.97547722973462
$ws.Cells.Item(19, 3).Value = 6.946804699397367
$ws.Cells.Item(19, 5).Value = 13.77326715588368
$ws.Cells.Item(19, 6).Value = 19.87204792380568
$ws.Cells.Item(19, 7).Value = 23.4186663146802
$ws.Cells.Item(19, 8).Value = 13.07121168654305
$ws.Cells.Item(19, 11).Value = 9.325416506949111
$ws.Cells.Item(19, 13).Value = 13.5584185366398
$ws.Cells.Item(19, 15).Value = 19.16219035656756
$ws.Cells.Item(20, 2).Value = 9.098637196264567
$ws.Cells.Item(20, 3).Value = 6.972527462686115
$ws.Cells.Item(20, 5).Value = 13.85997340338095
$ws.Cells.Item(20, 6).Value = 20.24955283636154
$ws.Cells.Item(20, 7).Value = 23.36505690844777
$ws.Cells.Item(20, 8).Value = 13.05106360049874
$ws.Cells.Item(20, 11).Value = 9.414865082194858
$ws.Cells.Item(20, 13).Value = 13.64002013472333
$ws.Cells.Item(20, 15).Value = 19.12537717218283
$ws.Cells.Item(21, 2).Value = 9.500874929187658
$ws.Cells.Item(21, 3).Value = 7.058406506099836
$ws.Cells.Item(21, 5).Value = 14.152822080892
$ws.Cells.Item(21, 6).Value = 21.46857628470577
$ws.Cells.Item(21, 7).Value = 23.19596521093012
$ws.Cells.Item(21, 8).Value = 12.98588357286457
$ws.Cells.Item(21, 11).Value = 9.70796507415306
$ws.Cells.Item(21, 13).Value = 13.91357966211969
$ws.Cells.Item(21, 15).Value = 19.00695062292164
$ws.Cells.Item(22, 2).Value = 9.754862808554691
$ws.Cells.Item(22, 3).Value = 7.114068633314158
$ws.Cells.Item(22, 5).Value = 14.34498619469631
$ws.Cells.Item(22, 6).Value = 22.22866616901552
$ws.Cells.Item(22, 7).Value = 23.09374547922972
$ws.Cells.Item(22, 8).Value = 12.94515716973854
$ws.Cells.Item(22, 11).Value = 9.893792638417809
$ws.Cells.Item(22, 13).Value = 14.0916149149844
$ws.Cells.Item(22, 15).Value = 18.93348063869568
$ws.Cells.Item(23, 2).Value = 9.620135417270486
$ws.Cells.Item(23, 3).Value = 7.084406917982807
$ws.Cells.Item(23, 5).Value = 14.24237739362315
$ws.Cells.Item(23, 6).Value = 21.82633154458858
$ws.Cells.Item(23, 7).Value = 23.14753878524421
$ws.Cells.Item(23, 8).Value = 12.96672239308146
$ws.Cells.Item(23, 11).Value = 9.795148765498192
$ws.Cells.Item(23, 13).Value = 13.99668163188978
$ws.Cells.Item(23, 15).Value = 18.97233297366541
$ws.Cells.Item(24, 2).Value = 9.091916567209706
$ws.Cells.Item(24, 3).Value = 6.971116796186097
$ws.Cells.Item(24, 5).Value = 13.85520521463294
$ws.Cells.Item(24, 6).Value = 20.22900810905287
$ws.Cells.Item(24, 7).Value = 23.36795957319389
$ws.Cells.Item(24, 8).Value = 13.05216057168485
$ws.Cells.Item(24, 11).Value = 9.409980456189357
$ws.Cells.Item(24, 13).Value = 13.63554049313139
$ws.Cells.Item(24, 15).Value = 19.12737900455648
$ws.Cells.Item(25, 2).Value = 8.485876911425159
$ws.Cells.Item(25, 3).Value = 6.847297117996688
$ws.Cells.Item(25, 5).Value = 13.44328881306046
$ws.Cells.Item(25, 6).Value = 18.34778573295695
$ws.Cells.Item(25, 7).Value = 23.64085394517465
$ws.Cells.Item(25, 8).Value = 13.15237101163763
$ws.Cells.Item(25, 11).Value = 8.971188152992282
$ws.Cells.Item(25, 13).Value = 13.24470915257766
$ws.Cells.Item(25, 15).Value = 19.31143948694925
